# Adds season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the rest of row 1 (bold, bordered, centered)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season record for every data row (2-58) with the team's record
$ws.Range("AD2:AD58").Value = 75
$ws.Range("AE2:AE58").Value = 87
$ws.Range("AF2:AF58").Value = 0
